$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New parameters: attaches pullrod sur inserts
$ws.Range("A15").Value = "Espacement_attaches_pullrod (mm)"
$ws.Range("B15").Value = 24.5

$ws.Range("A16").Value = "Rayon_attaches_pullrod (mm)"
$ws.Range("B16").Value = 2

# Match styling of the column above (text style col A, numeric style col B)
$ws.Range("A15:A16").NumberFormat = $ws.Range("A14").NumberFormat
$ws.Range("B15:B16").NumberFormat = $ws.Range("B14").NumberFormat

# Update active selection cell, as recorded in the saved workbook
$ws.Range("F9").Select()
